$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: extend merged title range C1:M1 -> C1:N1, add column N
$ws.Range("C1:M1").UnMerge()
$ws.Range("N1").Value2 = $ws.Range("M1").Value2
$ws.Range("C1:N1").Merge()
$ws.Range("C1:N1").Select()

# Row 2 (month labels) - new column N, label reused from original diff data ("дек" again)
$ws.Range("N2").Value2 = "дек"

# Row 3 - Потребительские цены (Consumer prices) - new January value
$ws.Range("N3").Value2 = "100.8"

# Row 4 - Продовольственные товары (Food products) - new January value
$ws.Range("N4").Value2 = "100.9"

# Row 5 - Непродовольственные товары (Non-food products) - new January value
$ws.Range("N5").Value2 = "100.5"

# Row 6 - Платные услуги (Paid services) - new January value
$ws.Range("N6").Value2 = "100.8"

# Row 7 - Цены производителей (Producer prices) - new January value
$ws.Range("N7").Value2 = "97.5"
